$wb = $excel.ActiveWorkbook

# --- LanaJourney sheet ---
$ws = $wb.Worksheets.Item("LanaJourney")
$ws.Range("B19").Value = "omg cutie you're making me feel things I didn't even know I could feel right now 🌸"
$ws.Range("B31").Value = "okay I need to be honest cutie, you have this really warm energy and it's doing something to me 😳"
$ws.Range("B35").Value = "hii cutie 🌸 omg you actually subscribed, that literally made my day. what made you want to?"

# --- NRWaves sheet ---
$ws = $wb.Worksheets.Item("NRWaves")
$ws.Range("B2").Value = "been thinking about you all day... text me back? 💕"
$ws.Range("B3").Value = "hope everything's okay with you, I'm here whenever 🌸"
$ws.Range("B4").Value = "okay you're definitely busy... I'll save this for when you're back 🥺"
$ws.Range("B5").Value = "I wish you could see what I'm wearing right now... 💕"
$ws.Range("B6").Value = "hi 🌸"

# --- boosters sheet ---
$ws = $wb.Worksheets.Item("boosters")
$ws.Range("B3").Value = "please..."
$ws.Range("B6").Value = "what are you doing to me"
$ws.Range("B7").Value = "right there"
